# Update "想去人数" (F column) figures on both the "展览" and "全部类型"
# sheets to reflect the latest scrape output.

$wb = $excel.ActiveWorkbook

# Sheet "展览" - exhibition-only listing
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 2415
$ws1.Range("F5").Value = 1781
$ws1.Range("F7").Value = 99
$ws1.Range("F8").Value = 838
$ws1.Range("F9").Value = 168

# Sheet "全部类型" - combined listing (same events, one row shifted down
# because it also includes a row from the "演出" sheet)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 2415
$ws4.Range("F5").Value = 1781
$ws4.Range("F8").Value = 99
$ws4.Range("F9").Value = 838
$ws4.Range("F10").Value = 168
